# Leave card update (1/9/2024 4:42 pm):
# Insert a new row above row 41 on Sheet1 (the leave-card ledger), splitting
# what used to be a single "UT(0-1-0)" entry (date 44713 / EARNED 1.25 /
# Absence-W-Pay 0.125) into two rows: the new row 41 keeps only the
# PARTICULARS label and the Absence-Undertime-W/Pay amount (date + EARNED
# cleared), while the original date/EARNED values now live on row 42 along
# with the same PARTICULARS/Absence value. Every row below shifts down by
# one (through the former last data row, which becomes row 133), and
# Table1 / the sheet dimension grow accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# Insert a blank row above row 41; this shifts rows 41:132 down to 42:133
# (values, formulas and styles all move down with the cells).
$ws.Rows.Item(41).Insert()

# The table's own range doesn't auto-grow from a plain row insert -- resize
# it to include the newly inserted row at the bottom of its old range.
$tbl.Resize($ws.Range("A8:K133"))

# Re-assert the calculated-column formula on the row that dropped out of
# the table range during the insert (row 133, formerly row 132) so it
# re-resolves its structured reference instead of caching #VALUE!.
$ws.Range("G133").Formula = $ws.Range("G133").Formula

# New row 41 copies its formatting from the row now below it (row 42, the
# shifted-down original row 41) since a bare row Insert() does not carry
# over number formats/borders in this runtime.
$ws.Range("A42:K42").Copy()
$ws.Range("A41:K41").PasteSpecial(-4122)   # xlPasteFormats

# Populate the new row 41 with its actual (edited) contents: only the
# PARTICULARS label and the Absence Undertime W/ Pay value remain; PERIOD
# and EARNED are blank.
$ws.Range("A41").Value2 = ""
$ws.Range("B41").Value2 = "UT(0-1-0)"
$ws.Range("C41").Value2 = ""
$ws.Range("D41").Value2 = 0.125
$ws.Range("G41").Formula = "=IF(ISBLANK([@EARNED]),"""",[@EARNED])"

# Restore the cursor/selection to where the edit left off.
$ws.Range("F40").Select() | Out-Null
